$wb = $excel.ActiveWorkbook

# --- Sheet 1: Evaluated_Interventions ---
$ws1 = $wb.Worksheets.Item("Evaluated_Interventions")

# Header row: metric label + year columns (years stored as text, matching source data)
$ws1.Range("A1").Value = "gta.evaluation"

$yearHeader1 = $ws1.Range("B1:K1")
$yearHeader1.NumberFormat = "@"
$ws1.Range("B1").Value = "2009"
$ws1.Range("C1").Value = "2010"
$ws1.Range("D1").Value = "2011"
$ws1.Range("E1").Value = "2012"
$ws1.Range("F1").Value = "2013"
$ws1.Range("G1").Value = "2014"
$ws1.Range("H1").Value = "2015"
$ws1.Range("I1").Value = "2016"
$ws1.Range("J1").Value = "2017"
$ws1.Range("K1").Value = "2018"
$yearHeader1.ClearFormats()

# Amber row
$ws1.Range("A2").Value = "Amber"
$ws1.Range("B2").Value = 191
$ws1.Range("C2").Value = 146
$ws1.Range("D2").Value = 142
$ws1.Range("E2").Value = 149
$ws1.Range("F2").Value = 160
$ws1.Range("G2").Value = 193
$ws1.Range("H2").Value = 226
$ws1.Range("I2").Value = 179
$ws1.Range("J2").Value = 160
$ws1.Range("K2").Value = 227

# Green row
$ws1.Range("A3").Value = "Green"
$ws1.Range("B3").Value = 391
$ws1.Range("C3").Value = 450
$ws1.Range("D3").Value = 469
$ws1.Range("E3").Value = 535
$ws1.Range("F3").Value = 490
$ws1.Range("G3").Value = 464
$ws1.Range("H3").Value = 521
$ws1.Range("I3").Value = 485
$ws1.Range("J3").Value = 448
$ws1.Range("K3").Value = 349

# Red row
$ws1.Range("A4").Value = "Red"
$ws1.Range("B4").Value = 1298
$ws1.Range("C4").Value = 1132
$ws1.Range("D4").Value = 1142
$ws1.Range("E4").Value = 1421
$ws1.Range("F4").Value = 1367
$ws1.Range("G4").Value = 1234
$ws1.Range("H4").Value = 1196
$ws1.Range("I4").Value = 1089
$ws1.Range("J4").Value = 1201
$ws1.Range("K4").Value = 1059

# --- Sheet 2: Total_Interventions ---
$ws2 = $wb.Worksheets.Item("Total_Interventions")

# Header row: year columns (years stored as text)
$yearHeader2 = $ws2.Range("A1:J1")
$yearHeader2.NumberFormat = "@"
$ws2.Range("A1").Value = "2009"
$ws2.Range("B1").Value = "2010"
$ws2.Range("C1").Value = "2011"
$ws2.Range("D1").Value = "2012"
$ws2.Range("E1").Value = "2013"
$ws2.Range("F1").Value = "2014"
$ws2.Range("G1").Value = "2015"
$ws2.Range("H1").Value = "2016"
$ws2.Range("I1").Value = "2017"
$ws2.Range("J1").Value = "2018"
$yearHeader2.ClearFormats()

# Totals row
$ws2.Range("A2").Value = 1880
$ws2.Range("B2").Value = 1728
$ws2.Range("C2").Value = 1753
$ws2.Range("D2").Value = 2105
$ws2.Range("E2").Value = 2017
$ws2.Range("F2").Value = 1891
$ws2.Range("G2").Value = 1943
$ws2.Range("H2").Value = 1753
$ws2.Range("I2").Value = 1809
$ws2.Range("J2").Value = 1635
